$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Local_User"
$ws.Range("A2").Value = "Net_172.16.0.0_24"
$ws.Range("B2").Value = "172.16.0.0/24"

$ws.Range("A3").Value = "Net_172.16.0.0_25"
$ws.Range("B3").Value = "172.16.0.0/25"
$ws.Range("C3").Value = "Local_User"

$ws.Range("C5").Select()
